$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for hour 20 and hour 21 (Lapaz column B)
$ws.Range("B21").Value = 23313
$ws.Range("B22").Value = 23000

# Add new rows for hour 22 (row 23) and hour 23 (row 24)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 23115
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 15564
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 22709
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 15389
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0

# Ensure the trailing empty row 25 is preserved in the output (matches
# the original <row r="25"/> placeholder element) even though rows 23/24
# now carry data.
$ws.Rows.Item(25).OutlineLevel = 0
